# Update the "Jogos da Semana" worksheet:
#  1. Add two new header columns (BC1/BD1) for HT correct-score odds 3-3 and 4-4,
#     matching the style used by the other header cells.
#  2. Refresh the single data row (row 2) with the latest match/odds data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers -----------------------------------------------------------
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# Copy the formatting of the existing last header cell (bold font, thin
# border, centered/top alignment) onto the two new header cells.
$ws.Range("BB1").Copy()
$ws.Range("BC1:BD1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2: identifying / descriptive fields --------------------------------
$ws.Range("A2").Value = "YyDJubM9"
$ws.Range("C2").Value = "11:30"
$ws.Range("D2").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E2").Value = "Al Qadisiya"
$ws.Range("F2").Value = "Al Feiha"

# --- Row 2: odds -------------------------------------------------------------
$ws.Range("G2").Value = 1.42
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 1.91
$ws.Range("K2").Value = 2.3
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.36
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 6.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 9.5
$ws.Range("AA2").Value = 13
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 12
$ws.Range("AD2").Value = 9
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 34
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 51
$ws.Range("AM2").Value = 51
$ws.Range("AN2").Value = 3.4
$ws.Range("AO2").Value = 7
$ws.Range("AP2").Value = 19
$ws.Range("AQ2").Value = 21
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 3
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 8
$ws.Range("AX2").Value = 34
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 151
$ws.Range("BB2").Value = 500
$ws.Range("BC2").Value = 81
$ws.Range("BD2").Value = 81
